$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.162.44"
$ws.Range("E2").Value = "  +2.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.419.07"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  +0.63%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.46"
$ws.Range("E5").Value = "  +1.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.14"
$ws.Range("E6").Value = "  +5.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.436.08"
$ws.Range("E9").Value = "  +1.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  +4.62%  "

$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("E12").Value = "  +4.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.349"
$ws.Range("E13").Value = "  +4.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.65"
$ws.Range("E14").Value = "  +4.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000177"
$ws.Range("E15").Value = "  +6.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.845.18"
$ws.Range("E16").Value = "  +1.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.042.29"
$ws.Range("E17").Value = "  +2.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.440.40"
$ws.Range("E18").Value = "  +2.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.91"
$ws.Range("E19").Value = "  -2.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.82"
$ws.Range("E20").Value = "  +2.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.91"
$ws.Range("E21").Value = "  +1.25%  "

$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.08"
$ws.Range("E23").Value = "  -1.00%  "

$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  +13.72%  "

$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.41"
$ws.Range("E26").Value = "  +2.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "618.20"
$ws.Range("E27").Value = "  +11.51%  "

$ws.Range("E28").Value = "  +4.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0989"
$ws.Range("E29").Value = "  +8.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.537.86"
$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.05"
$ws.Range("E31").Value = "  +2.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("E32").Value = "  +10.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.82"
$ws.Range("E33").Value = "  +1.77%  "

$ws.Range("E34").Value = "  +3.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.48"
$ws.Range("E35").Value = "  +5.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  -0.67%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.00"
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.64"
$ws.Range("E38").Value = "  +3.30%  "

$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.372"
$ws.Range("E39").Value = "  +1.38%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.37"
$ws.Range("E40").Value = "  +6.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.51"
$ws.Range("E41").Value = "  +2.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.66"
$ws.Range("E42").Value = "  +16.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.72"
$ws.Range("E43").Value = "  +5.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.30"
$ws.Range("E44").Value = "  +3.05%  "

$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0284"
$ws.Range("E46").Value = "  -1.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.19"
$ws.Range("E47").Value = "  +0.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.57"
$ws.Range("E48").Value = "  +2.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.26"
$ws.Range("E49").Value = "  +7.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.599"
$ws.Range("E50").Value = "  +2.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0513"
$ws.Range("E51").Value = "  +3.51%  "
